$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap full row content (columns B:AC); column A (match id) stays in place ---
# Swap row 16 <-> row 17
$ws.Range("B16").Value = 7138607
$ws.Range("C16").Value = 'Germany Verbandsliga'
$ws.Range("D16").Value = 'Germany Verbandsliga'
$ws.Range("E16").Value = 45168.625
$ws.Range("F16").Value = 'Rot Weiss Walldorf II'
$ws.Range("G16").Value = 'Turnerschaft OberRoden'
$ws.Range("H16").Value = 3
$ws.Range("I16").Value = 2
$ws.Range("J16").Value = 'H'
$ws.Range("K16").Value = 2.25
$ws.Range("L16").Value = 3.75
$ws.Range("M16").Value = 2.5
$ws.Range("N16").Value = 2.25
$ws.Range("O16").Value = 3.8
$ws.Range("P16").Value = 2.45
$ws.Range("Q16").Value = 0
$ws.Range("R16").Value = 1.8
$ws.Range("S16").Value = 2
$ws.Range("T16").Value = 3.75
$ws.Range("U16").Value = 1.95
$ws.Range("V16").Value = 1.85
$ws.Range("W16").Value = 1.25
$ws.Range("X16").Value = -1
$ws.Range("Y16").Value = -1
$ws.Range("Z16").Value = 0.8
$ws.Range("AA16").Value = -1
$ws.Range("AB16").Value = 0.95
$ws.Range("AC16").Value = -1
$ws.Range("B17").Value = 7138608
$ws.Range("C17").Value = 'Germany Verbandsliga'
$ws.Range("D17").Value = 'Germany Verbandsliga'
$ws.Range("E17").Value = 45168.625
$ws.Range("F17").Value = 'SV UnterFlockenbach'
$ws.Range("G17").Value = 'SC Dortelweil'
$ws.Range("H17").Value = 1
$ws.Range("I17").Value = 1
$ws.Range("J17").Value = 'D'
$ws.Range("K17").Value = 1.083
$ws.Range("L17").Value = 9
$ws.Range("M17").Value = 16
$ws.Range("N17").Value = 1.125
$ws.Range("O17").Value = 7.5
$ws.Range("P17").Value = 13
$ws.Range("Q17").Value = -2.5
$ws.Range("R17").Value = 1.775
$ws.Range("S17").Value = 1.925
$ws.Range("T17").Value = 4.25
$ws.Range("U17").Value = 1.975
$ws.Range("V17").Value = 1.825
$ws.Range("W17").Value = -1
$ws.Range("X17").Value = 6.5
$ws.Range("Y17").Value = -1
$ws.Range("Z17").Value = -1
$ws.Range("AA17").Value = 0.925
$ws.Range("AB17").Value = -1
$ws.Range("AC17").Value = 0.825

# Swap row 20 <-> row 21
$ws.Range("B20").Value = 7149166
$ws.Range("C20").Value = 'Germany Verbandsliga'
$ws.Range("D20").Value = 'Germany Verbandsliga'
$ws.Range("E20").Value = 45170.60416666666
$ws.Range("F20").Value = 'Fuchse Berlin Reinickendorf'
$ws.Range("G20").Value = 'SD Croatia Berlin'
$ws.Range("H20").Value = 1
$ws.Range("I20").Value = 4
$ws.Range("J20").Value = 'A'
$ws.Range("K20").Value = 1.4
$ws.Range("L20").Value = 4.8
$ws.Range("M20").Value = 5.25
$ws.Range("N20").Value = 1.4
$ws.Range("O20").Value = 5
$ws.Range("P20").Value = 5
$ws.Range("Q20").Value = -1.25
$ws.Range("R20").Value = 1.8
$ws.Range("S20").Value = 2
$ws.Range("T20").Value = 3.5
$ws.Range("U20").Value = 1.925
$ws.Range("V20").Value = 1.875
$ws.Range("W20").Value = -1
$ws.Range("X20").Value = -1
$ws.Range("Y20").Value = 4
$ws.Range("Z20").Value = -1
$ws.Range("AA20").Value = 1
$ws.Range("AB20").Value = 0.925
$ws.Range("AC20").Value = -1
$ws.Range("B21").Value = 7149361
$ws.Range("C21").Value = 'Germany Verbandsliga'
$ws.Range("D21").Value = 'Germany Verbandsliga'
$ws.Range("E21").Value = 45170.60416666666
$ws.Range("F21").Value = 'TSG 1846 Bretzenheim'
$ws.Range("G21").Value = 'TSV 1881 GauOdernheim'
$ws.Range("H21").Value = 4
$ws.Range("I21").Value = 2
$ws.Range("J21").Value = 'H'
$ws.Range("K21").Value = 3.75
$ws.Range("L21").Value = 4.333
$ws.Range("M21").Value = 1.615
$ws.Range("N21").Value = 3.75
$ws.Range("O21").Value = 4.333
$ws.Range("P21").Value = 1.615
$ws.Range("Q21").Value = 0.75
$ws.Range("R21").Value = 2
$ws.Range("S21").Value = 1.8
$ws.Range("T21").Value = 3.75
$ws.Range("U21").Value = 1.9
$ws.Range("V21").Value = 1.9
$ws.Range("W21").Value = 2.75
$ws.Range("X21").Value = -1
$ws.Range("Y21").Value = -1
$ws.Range("Z21").Value = 1
$ws.Range("AA21").Value = -1
$ws.Range("AB21").Value = 0.8999999999999999
$ws.Range("AC21").Value = -1

# Swap row 67 <-> row 68
$ws.Range("B67").Value = 7423702
$ws.Range("C67").Value = 'Germany Verbandsliga'
$ws.Range("D67").Value = 'Germany Verbandsliga'
$ws.Range("E67").Value = 45233.625
$ws.Range("F67").Value = 'SG RotWeiss Thalheim'
$ws.Range("G67").Value = 'SV Fortuna Magdeburg'
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 'D'
$ws.Range("K67").Value = 1.8
$ws.Range("L67").Value = 4.5
$ws.Range("M67").Value = 3
$ws.Range("N67").Value = 1.8
$ws.Range("O67").Value = 4.5
$ws.Range("P67").Value = 3
$ws.Range("Q67").Value = -0.5
$ws.Range("R67").Value = 1.85
$ws.Range("S67").Value = 1.95
$ws.Range("T67").Value = 3.5
$ws.Range("U67").Value = 1.85
$ws.Range("V67").Value = 1.95
$ws.Range("W67").Value = -1
$ws.Range("X67").Value = 3.5
$ws.Range("Y67").Value = -1
$ws.Range("Z67").Value = -1
$ws.Range("AA67").Value = 0.95
$ws.Range("AB67").Value = -1
$ws.Range("AC67").Value = 0.95
$ws.Range("B68").Value = 7423701
$ws.Range("C68").Value = 'Germany Verbandsliga'
$ws.Range("D68").Value = 'Germany Verbandsliga'
$ws.Range("E68").Value = 45233.625
$ws.Range("F68").Value = 'BSV HalleAmmendorf'
$ws.Range("G68").Value = 'VfB Sangerhausen'
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 2
$ws.Range("J68").Value = 'A'
$ws.Range("K68").Value = 2
$ws.Range("L68").Value = 4.5
$ws.Range("M68").Value = 2.55
$ws.Range("N68").Value = 2
$ws.Range("O68").Value = 4.5
$ws.Range("P68").Value = 2.6
$ws.Range("Q68").Value = -0.25
$ws.Range("R68").Value = 1.85
$ws.Range("S68").Value = 1.95
$ws.Range("T68").Value = 3.25
$ws.Range("U68").Value = 1.85
$ws.Range("V68").Value = 1.95
$ws.Range("W68").Value = -1
$ws.Range("X68").Value = -1
$ws.Range("Y68").Value = 1.6
$ws.Range("Z68").Value = -1
$ws.Range("AA68").Value = 0.95
$ws.Range("AB68").Value = -1
$ws.Range("AC68").Value = 0.95

# --- Single-cell team-name updates (shared-string re-sort side effects) ---
$ws.Range("G2").Value = 'SV Frankonia Wernsdorf'
$ws.Range("G3").Value = 'SG RotWeiss Thalheim'
$ws.Range("G4").Value = '1 FC BitterfeldWolfen'
$ws.Range("G5").Value = 'SV 1908 GW Ahrensfelde'
$ws.Range("G6").Value = 'SV Frankonia Wernsdorf'
$ws.Range("G7").Value = 'MSC Preussen 1899'
$ws.Range("G8").Value = 'TSV Mariendorf 1897'
$ws.Range("G9").Value = 'Berlin Turkspor'
$ws.Range("G10").Value = 'VSG Altglienicke II'
$ws.Range("G15").Value = 'JSK Rodgau'
$ws.Range("G41").Value = 'Berlin Turkspor'
$ws.Range("F95").Value = 'SG RotWeiss Thalheim'
$ws.Range("F98").Value = 'BSV HalleAmmendorf'

# --- Append new match rows 101-103 ---
# Row 101: copy formatting from row 100, then set values
$ws.Range("A100:AC100").Copy() | Out-Null
$ws.Range("A101:AC101").PasteSpecial(-4122) | Out-Null
# Row 102: copy formatting from row 101, then set values
$ws.Range("A101:AC101").Copy() | Out-Null
$ws.Range("A102:AC102").PasteSpecial(-4122) | Out-Null
# Row 103: copy formatting from row 102, then set values
$ws.Range("A102:AC102").Copy() | Out-Null
$ws.Range("A103:AC103").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A101").Value = 99
$ws.Range("B101").Value = 7905680
$ws.Range("C101").Value = 'Germany Verbandsliga'
$ws.Range("D101").Value = 'Germany Verbandsliga'
$ws.Range("E101").Value = 45354.45833333334
$ws.Range("F101").Value = 'RotWeiss Frankfurt'
$ws.Range("G101").Value = 'SV Pars NeuIsenburg'
$ws.Range("H101").Value = 1
$ws.Range("I101").Value = 2
$ws.Range("J101").Value = 'A'
$ws.Range("K101").Value = 3.5
$ws.Range("L101").Value = 4
$ws.Range("M101").Value = 1.727
$ws.Range("N101").Value = 3.5
$ws.Range("O101").Value = 4
$ws.Range("P101").Value = 1.727
$ws.Range("Q101").Value = 0.75
$ws.Range("R101").Value = 1.85
$ws.Range("S101").Value = 1.95
$ws.Range("T101").Value = 4
$ws.Range("U101").Value = 1.975
$ws.Range("V101").Value = 1.825
$ws.Range("W101").Value = -1
$ws.Range("X101").Value = -1
$ws.Range("Y101").Value = 0.7270000000000001
$ws.Range("Z101").Value = -0.5
$ws.Range("AA101").Value = 0.475
$ws.Range("AB101").Value = -1
$ws.Range("AC101").Value = 0.825

$ws.Range("A102").Value = 100
$ws.Range("B102").Value = 7905679
$ws.Range("C102").Value = 'Germany Verbandsliga'
$ws.Range("D102").Value = 'Germany Verbandsliga'
$ws.Range("E102").Value = 45354.47916666666
$ws.Range("F102").Value = 'FFV Sportfreunde 04'
$ws.Range("G102").Value = 'DJK Bad Homburg'
$ws.Range("H102").Value = 1
$ws.Range("I102").Value = 5
$ws.Range("J102").Value = 'A'
$ws.Range("K102").Value = 7
$ws.Range("L102").Value = 6
$ws.Range("M102").Value = 1.25
$ws.Range("N102").Value = 7
$ws.Range("O102").Value = 6
$ws.Range("P102").Value = 1.25
$ws.Range("Q102").Value = 2
$ws.Range("R102").Value = 1.85
$ws.Range("S102").Value = 1.95
$ws.Range("T102").Value = 4
$ws.Range("U102").Value = 1.95
$ws.Range("V102").Value = 1.85
$ws.Range("W102").Value = -1
$ws.Range("X102").Value = -1
$ws.Range("Y102").Value = 0.25
$ws.Range("Z102").Value = -1
$ws.Range("AA102").Value = 0.95
$ws.Range("AB102").Value = 0.95
$ws.Range("AC102").Value = -1

$ws.Range("A103").Value = 101
$ws.Range("B103").Value = 7919896
$ws.Range("C103").Value = 'Germany Verbandsliga'
$ws.Range("D103").Value = 'Germany Verbandsliga'
$ws.Range("E103").Value = 45357.66666666666
$ws.Range("F103").Value = 'JSK Rodgau'
$ws.Range("G103").Value = 'RotWeiss Darmstadt'
$ws.Range("H103").Value = 1
$ws.Range("I103").Value = 2
$ws.Range("J103").Value = 'A'
$ws.Range("K103").Value = 2.75
$ws.Range("L103").Value = 4
$ws.Range("M103").Value = 2
$ws.Range("N103").Value = 3.75
$ws.Range("O103").Value = 4.333
$ws.Range("P103").Value = 1.615
$ws.Range("Q103").Value = 0.75
$ws.Range("R103").Value = 2
$ws.Range("S103").Value = 1.8
$ws.Range("T103").Value = 4
$ws.Range("U103").Value = 1.825
$ws.Range("V103").Value = 1.975
$ws.Range("W103").Value = -1
$ws.Range("X103").Value = -1
$ws.Range("Y103").Value = 0.615
$ws.Range("Z103").Value = -0.5
$ws.Range("AA103").Value = 0.4
$ws.Range("AB103").Value = -1
$ws.Range("AC103").Value = 0.9750000000000001

